# Append: 2025-09-19 12:45 JST
# Update the "取得日時" (acquired timestamp) column for the existing
# ランサーズ rows from the previous run's timestamp to the
# newly scraped timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-19 12:45:57"

for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
